# Add new query criterion "user.decryptUntrusted" to the criterion_property sheet.
# Alphabetically it belongs right after "user.decrypt" (row 269) and before
# "user.deferredDelete" (previously row 270) - so it becomes the new row 270,
# pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("criterion_property")
$ws.Activate()

# Insert a fresh blank row at position 270; old row 270 ("user.deferredDelete")
# and everything after it shifts down to row 271 onward.
$ws.Rows.Item(270).Insert()
$ws.Rows.Item(270).RowHeight = 16.5

# Populate the new row 270 with the new criterion, following the same simple
# BOOLEAN-typed pattern used by its neighbours "user.decrypt" / "user.deferredDelete":
#   A = #module, B = property, C = value_type, K = name_l10n_key, L = valid_restrictions
$ws.Range("A270").Value = "USER_DB"
$ws.Range("B270").Value = "user.decryptUntrusted"
$ws.Range("C270").Value = "BOOLEAN"
$ws.Range("K270").Value = "user.decryptUntrusted"
$ws.Range("L270").Value = "EQ, NE"

# Reflect the edit location in the view/selection.
$excel.ActiveWindow.ScrollRow = 248
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A270").Select()
